$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 288, shifting existing rows 288-320 down to 289-321.
$ws.Rows.Item(288).Insert()

# Populate the newly inserted row 288 with the new record's data.
$ws.Cells.Item(288, 1).Value = 5
$ws.Cells.Item(288, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(288, 3).Value = "Maule"
$ws.Cells.Item(288, 4).Value = 44918
$ws.Cells.Item(288, 5).Value = 7
$ws.Cells.Item(288, 6).Value = "Fruta"
$ws.Cells.Item(288, 7).Value = 100108
$ws.Cells.Item(288, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(288, 9).Value = 100108005
$ws.Cells.Item(288, 10).Value = "Piña"
$ws.Cells.Item(288, 11).Value = "Caramelo"
$ws.Cells.Item(288, 12).Value = "Segunda"
$ws.Cells.Item(288, 13).Value = 230
$ws.Cells.Item(288, 14).Value = 19000
$ws.Cells.Item(288, 15).Value = 19000
$ws.Cells.Item(288, 16).Value = 19000
$ws.Cells.Item(288, 17).Value = "`$/caja 14 unidades"
$ws.Cells.Item(288, 18).Value = "Ecuador"
$ws.Cells.Item(288, 19).Value = 1357
$ws.Cells.Item(288, 20).Value = 14
